# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D) for the
# 055ae8aa-8ee6-4a68-a56c-53ccce15ef9c file row (row 4) on both the
# zh-cn and de-de status sheets to reflect the new handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-03-11 06:49:51"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-03-11 06:50:01"
